$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu results for Case_4_88 (380 kV case), rows 2-25 (bus 0-23)
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029358489464243
$ws.Range("D2").Value = 1.039328387759727
$ws.Range("E2").Value = 1.032993435636037
$ws.Range("F2").Value = 1.047802589758732
$ws.Range("I2").Value = 1.036636729298295
$ws.Range("J2").Value = 1.034505705879688
$ws.Range("K2").Value = 1.042113887532002
$ws.Range("L2").Value = 1.035797065510537
$ws.Range("M2").Value = 1.050564208927227
$ws.Range("N2").Value = 1.01549454372717

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030314250437209
$ws.Range("D3").Value = 1.04008023098506
$ws.Range("E3").Value = 1.033896929726259
$ws.Range("F3").Value = 1.048687320335965
$ws.Range("I3").Value = 1.036841210322978
$ws.Range("J3").Value = 1.035102433099425
$ws.Range("K3").Value = 1.042676080741799
$ws.Range("L3").Value = 1.036509194597152
$ws.Range("M3").Value = 1.051260668364345
$ws.Range("N3").Value = 1.015693505177421

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030933013231556
$ws.Range("D4").Value = 1.040566633895493
$ws.Range("E4").Value = 1.034482223325876
$ws.Range("F4").Value = 1.049260032161022
$ws.Range("I4").Value = 1.036971835583761
$ws.Range("J4").Value = 1.035488277858751
$ws.Range("K4").Value = 1.043039104636248
$ws.Range("L4").Value = 1.036970032394771
$ws.Range("M4").Value = 1.051710919707445
$ws.Range("N4").Value = 1.015822107183003

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031193217281273
$ws.Range("D5").Value = 1.040771094823307
$ws.Range("E5").Value = 1.034728440573235
$ws.Range("F5").Value = 1.04950085440024
$ws.Range("I5").Value = 1.037026345758757
$ws.Range("J5").Value = 1.035650419654175
$ws.Range("K5").Value = 1.043191538654886
$ws.Range("L5").Value = 1.037163777823753
$ws.Range("M5").Value = 1.051900107373743
$ws.Range("N5").Value = 1.015876137767179

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031236911129507
$ws.Range("D6").Value = 1.04080542333384
$ws.Range("E6").Value = 1.034769790869901
$ws.Range("F6").Value = 1.049541292652195
$ws.Range("I6").Value = 1.037035474506811
$ws.Range("J6").Value = 1.035677640000812
$ws.Range("K6").Value = 1.043217122340198
$ws.Range("L6").Value = 1.03719630904547
$ws.Range("M6").Value = 1.051931867034793
$ws.Range("N6").Value = 1.015885207754083

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030936489792509
$ws.Range("D7").Value = 1.0405693660039
$ws.Range("E7").Value = 1.034485512666296
$ws.Range("F7").Value = 1.049263249828349
$ws.Range("I7").Value = 1.036972565542356
$ws.Range("J7").Value = 1.035490444671059
$ws.Range("K7").Value = 1.04304114218046
$ws.Range("L7").Value = 1.036972621195704
$ws.Range("M7").Value = 1.051713448028839
$ws.Range("N7").Value = 1.015822829275152

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029681426299789
$ws.Range("D8").Value = 1.039582494522878
$ws.Range("E8").Value = 1.03329863573683
$ws.Range("F8").Value = 1.048101539756717
$ws.Range("I8").Value = 1.036706183657443
$ws.Range("J8").Value = 1.034707429604893
$ws.Range("K8").Value = 1.042304038595554
$ws.Range("L8").Value = 1.03603772343962
$ws.Range("M8").Value = 1.050799663707234
$ws.Range("N8").Value = 1.015561812333478

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02747234195202
$ws.Range("D9").Value = 1.037842870344375
$ws.Range("E9").Value = 1.031212409506195
$ws.Range("F9").Value = 1.046056294113215
$ws.Range("I9").Value = 1.036223888434456
$ws.Range("J9").Value = 1.033325568722247
$ws.Range("K9").Value = 1.040999448138313
$ws.Range("L9").Value = 1.034390684524901
$ws.Range("M9").Value = 1.049186412200731
$ws.Range("N9").Value = 1.015100814612007

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026001344234116
$ws.Range("D10").Value = 1.03668277421624
$ws.Range("E10").Value = 1.029825158647949
$ws.Range("F10").Value = 1.044694109829056
$ws.Range("I10").Value = 1.03589372579592
$ws.Range("J10").Value = 1.032402974459701
$ws.Range("K10").Value = 1.040125935091241
$ws.Range("L10").Value = 1.033292965251719
$ws.Range("M10").Value = 1.048108930279531
$ws.Range("N10").Value = 1.014792794721361

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025364804217437
$ws.Range("D11").Value = 1.036180373464728
$ws.Range("E11").Value = 1.029225323742516
$ws.Range("F11").Value = 1.044104596470485
$ws.Range("I11").Value = 1.035748722017881
$ws.Range("J11").Value = 1.032003170144114
$ws.Range("K11").Value = 1.039746808474387
$ws.Range("L11").Value = 1.032817725195156
$ws.Range("M11").Value = 1.047641912698147
$ws.Range("N11").Value = 1.014659259519732

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02512842742354
$ws.Range("D12").Value = 1.035993749955104
$ws.Range("E12").Value = 1.029002647725759
$ws.Range("F12").Value = 1.04388567449855
$ws.Range("I12").Value = 1.035694555075162
$ws.Range("J12").Value = 1.031854618488271
$ws.Range("K12").Value = 1.039605851230595
$ws.Range("L12").Value = 1.032641212744134
$ws.Range("M12").Value = 1.047468373406819
$ws.Range("N12").Value = 1.014609634786785

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025179128208953
$ws.Range("D13").Value = 1.036033781741565
$ws.Range("E13").Value = 1.029050406615276
$ws.Range("F13").Value = 1.043932631743021
$ws.Range("I13").Value = 1.035706187909867
$ws.Range("J13").Value = 1.031886485415871
$ws.Range("K13").Value = 1.039636093029641
$ws.Range("L13").Value = 1.032679074667381
$ws.Range("M13").Value = 1.047505601241615
$ws.Range("N13").Value = 1.014620280533121

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02534526395892
$ws.Range("D14").Value = 1.036164947293889
$ws.Range("E14").Value = 1.029206914628
$ws.Range("F14").Value = 1.044086499294244
$ws.Range("I14").Value = 1.035744250802155
$ws.Range("J14").Value = 1.031990891762628
$ws.Range("K14").Value = 1.039735159607084
$ws.Range("L14").Value = 1.032803134353215
$ws.Range("M14").Value = 1.047627569270659
$ws.Range("N14").Value = 1.014655158007832

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025447633980536
$ws.Range("D15").Value = 1.03624576149761
$ws.Range("E15").Value = 1.029303361556057
$ws.Range("F15").Value = 1.044181308772261
$ws.Range("I15").Value = 1.03576766206465
$ws.Range("J15").Value = 1.032055213816112
$ws.Range("K15").Value = 1.039796180233792
$ws.Range("L15").Value = 1.03287957335198
$ws.Range("M15").Value = 1.047702708797646
$ws.Range("N15").Value = 1.014676644024807

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026043598028732
$ws.Range("D16").Value = 1.036716115499339
$ws.Range("E16").Value = 1.029864985798296
$ws.Range("F16").Value = 1.044733240767741
$ws.Range("I16").Value = 1.035903306263463
$ws.Range("J16").Value = 1.032429501587207
$ws.Range("K16").Value = 1.04015107781063
$ws.Range("L16").Value = 1.033324507123919
$ws.Range("M16").Value = 1.048139915080564
$ws.Range("N16").Value = 1.014801653656507

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026417541218554
$ws.Range("D17").Value = 1.037011137877522
$ws.Range("E17").Value = 1.030217507223309
$ws.Range("F17").Value = 1.045079540062981
$ws.Range("I17").Value = 1.035987846215316
$ws.Range("J17").Value = 1.032664198658034
$ws.Range("K17").Value = 1.040373458143342
$ws.Range("L17").Value = 1.033603624281449
$ws.Range("M17").Value = 1.04841404041307
$ws.Range("N17").Value = 1.014880026128232

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026635695633757
$ws.Range("D18").Value = 1.037183212460602
$ws.Range("E18").Value = 1.030423209319454
$ws.Range("F18").Value = 1.045281561512478
$ws.Range("I18").Value = 1.036036959956177
$ws.Range("J18").Value = 1.032801063052059
$ws.Range("K18").Value = 1.040503082923406
$ws.Range("L18").Value = 1.033766436220141
$ws.Range("M18").Value = 1.048573888550221
$ws.Range("N18").Value = 1.014925723953039

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02671008736205
$ws.Range("D19").Value = 1.037241884225211
$ws.Range("E19").Value = 1.030493362346955
$ws.Range("F19").Value = 1.04535045081424
$ws.Range("I19").Value = 1.036053673043209
$ws.Range("J19").Value = 1.032847725087186
$ws.Range("K19").Value = 1.040547267024167
$ws.Range("L19").Value = 1.033821952149317
$ws.Range("M19").Value = 1.048628385048464
$ws.Range("N19").Value = 1.014941303089864

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026377416544232
$ws.Range("D20").Value = 1.036979485481267
$ws.Range("E20").Value = 1.030179676488677
$ws.Range("F20").Value = 1.045042382217699
$ws.Range("I20").Value = 1.035978796249185
$ws.Range("J20").Value = 1.032639021026009
$ws.Range("K20").Value = 1.040349607717084
$ws.Range("L20").Value = 1.03357367685694
$ws.Range("M20").Value = 1.048384633957561
$ws.Range("N20").Value = 1.014871619105963

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02529633939051
$ws.Range("D21").Value = 1.036126322560313
$ws.Range("E21").Value = 1.029160823338988
$ws.Range("F21").Value = 1.044041187757984
$ws.Range("I21").Value = 1.035733050673816
$ws.Range("J21").Value = 1.031960147974015
$ws.Range("K21").Value = 1.039705990625954
$ws.Range("L21").Value = 1.032766601505442
$ws.Range("M21").Value = 1.047591654599327
$ws.Range("N21").Value = 1.014644888115722

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024616985267517
$ws.Range("D22").Value = 1.035589850808311
$ws.Range("E22").Value = 1.028520978986756
$ws.Range("F22").Value = 1.043411984287311
$ws.Range("I22").Value = 1.035576770059284
$ws.Range("J22").Value = 1.031533044708561
$ws.Range("K22").Value = 1.039300555494404
$ws.Range("L22").Value = 1.03225923608368
$ws.Range("M22").Value = 1.047092683129897
$ws.Range("N22").Value = 1.014502195718299

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024977089209628
$ws.Range("D23").Value = 1.035874249375548
$ws.Range("E23").Value = 1.028860101218146
$ws.Range("F23").Value = 1.043745509311191
$ws.Range("I23").Value = 1.035659785010874
$ws.Range("J23").Value = 1.031759485452861
$ws.Range("K23").Value = 1.039515556771198
$ws.Range("L23").Value = 1.032528192694446
$ws.Range("M23").Value = 1.047357234327618
$ws.Range("N23").Value = 1.014577852600495

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026395547034801
$ws.Range("D24").Value = 1.036993787858334
$ws.Range("E24").Value = 1.030196770314995
$ws.Range("F24").Value = 1.045059172153064
$ws.Range("I24").Value = 1.035982886148003
$ws.Range("J24").Value = 1.032650397808249
$ws.Range("K24").Value = 1.040360384963572
$ws.Range("L24").Value = 1.033587208785504
$ws.Range("M24").Value = 1.048397921606683
$ws.Range("N24").Value = 1.014875417925401

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028043142484657
$ws.Range("D25").Value = 1.03829267163404
$ws.Range("E25").Value = 1.031751125373817
$ws.Range("F25").Value = 1.046584813470473
$ws.Range("I25").Value = 1.036350097635585
$ws.Range("J25").Value = 1.033683055157635
$ws.Range("K25").Value = 1.041337387710108
$ws.Range("L25").Value = 1.034816433723312
$ws.Range("M25").Value = 1.049603830465733
$ws.Range("N25").Value = 1.015220116281212
